$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update trieda for Tomas (row 5): I.A -> 1.A
$ws.Range("D5").Value = "1.A"

# Rename student in row 6 from Betka to Lumi
$ws.Range("B6").Value = "Ľumi"

# Update trieda for Judas (row 7): Prima with diacritics -> without
$ws.Range("D7").Value = "Prima"

# Update trieda (class) for Maros (row 4): Oktava with diacritics -> without
$ws.Range("D4").Value = "Oktava"

# Update trieda for Hanna (row 8): Oktava with diacritics -> without
$ws.Range("D8").Value = "Oktava"

# Remove the last student row (Samo, id 10)
$ws.Rows(11).Delete() | Out-Null

# Update the current selection to match the new workbook state
$ws.Range("M17:Q24").Select() | Out-Null
